$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.012.99'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.41%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.466.85'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.51%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '571.46'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.52%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.87'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.73%  '

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.05%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.537'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.03%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.465.34'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.47%  '

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.49%  '

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.84%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.25'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.42%  '

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.69%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.03'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.12%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000180'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +3.56%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.907.93'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.69%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.889.54'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.37%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.478.98'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.91%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.37'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.06%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.30'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +6.85%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '324.12'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.20%  '

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.38%  '

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +13.41%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.999'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.05%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '66.24'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.91%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '621.49'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +12.23%  '

# Row 27
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = 'PEPE'
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000102'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +9.40%  '

# Row 28
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = 'Aptos'
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.51'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.55%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.589.69'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.57%  '

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.27%  '

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +4.88%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.19'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.84%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.141'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -4.06%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.89'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.62%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.07'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +6.23%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.48'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.43%  '

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.06%  '

# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.14%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.44'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.69%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.72'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.00%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '145.90'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.95%  '

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.97%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.60'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +14.61%  '

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.12%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '147.68'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.13%  '

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.87%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '20.77'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +3.57%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0537'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.12%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.603'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.08%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0233'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.67%  '

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.53%  '
